$wb = $excel.ActiveWorkbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws.Name = "grupos_has_usuarios(MySQL)"

# Write the "Y" (principal) rows first, in ascending group order, so the
# shared-string table gets these 20 strings appended before the "N" rows
# (matches the order they were authored in upstream).
$yRows = @(1,6,11,16,21,26,31,36,41,46,51,56,61,66,71,76,81,86,91,96)
$yValues = @(
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (26,1,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (97,2,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (6,3,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (22,4,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (24,5,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (58,6,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (80,7,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (94,8,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (23,9,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (7,10,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (14,11,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (51,12,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (9,13,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (5,14,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (5,15,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (21,16,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (18,17,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (56,18,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (96,19,"Y");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (44,20,"Y");'
)
for ($i = 0; $i -lt $yRows.Length; $i++) {
    $ws.Cells.Item($yRows[$i], 1).Value = $yValues[$i]
}

# Then write the "N" rows, in ascending group order.
$nRows = @(2,3,4,5,7,8,9,10,12,13,14,15,17,18,19,20,22,23,24,25,27,28,29,30,32,33,34,35,37,38,39,40,42,43,44,45,47,48,49,50,52,53,54,55,57,58,59,60,62,63,64,65,67,68,69,70,72,73,74,75,77,78,79,80,82,83,84,85,87,88,89,90,92,93,94,95,97,98,99,100)
$nValues = @(
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (82,1,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (81,1,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (99,1,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (65,1,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (3,2,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (42,2,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (58,2,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (80,2,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (28,3,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (4,3,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (92,3,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (19,3,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (42,4,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (64,4,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (11,4,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (74,4,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (96,5,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (85,5,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (48,5,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (72,5,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (19,6,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (32,6,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (33,6,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (67,6,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (9,7,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (67,7,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (51,7,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (83,7,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (51,8,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (70,8,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (84,8,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (48,8,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (46,9,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (24,9,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (35,9,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (13,9,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (92,10,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (75,10,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (9,10,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (42,10,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (23,11,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (78,11,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (8,11,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (19,11,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (36,12,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (77,12,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (16,12,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (29,12,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (38,13,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (70,13,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (100,13,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (23,13,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (43,14,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (12,14,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (56,14,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (23,14,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (64,15,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (1,15,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (31,15,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (46,15,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (62,16,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (71,16,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (98,16,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (99,16,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (74,17,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (97,17,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (83,17,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (62,17,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (98,18,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (56,18,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (38,18,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (13,18,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (43,19,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (65,19,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (56,19,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (98,19,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (57,20,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (44,20,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (88,20,"N");',
    'INSERT INTO `grupos_has_usuarios` (`grupos_idgrupo`,`usuarios_idusuario`,`principal`) VALUES (79,20,"N");'
)
for ($i = 0; $i -lt $nRows.Length; $i++) {
    $ws.Cells.Item($nRows[$i], 1).Value = $nValues[$i]
}

Write-Output "Y rows: $($yRows.Length), N rows: $($nRows.Length)"
